$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold numeric-looking text (e.g. "0.620", "226.92").
# Force a Text number format first so Excel keeps the exact original
# string (incl. trailing zeros) instead of auto-converting to a number.
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D14", "D15", "D16", "D19", "D21", "D22", "D25", "D26", "D28", "D29", "D30", "D32", "D33", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D46", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '39.522.00'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '2.161.56'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = '226.92'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").Value = '0.620'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '62.57'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").Value = '0.0845'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").Value = '15.81'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = '2.483.88'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Value = '21.64'
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").Value = '0.804'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '5.44'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '2.163.58'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '39.537.26'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '71.53'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '0.0₃0886'
$ws.Range("E20").Value = '  +4.23%  '
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  -2.56%  '
$ws.Range("D22").Value = '227.32'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").Value = '2.32'
$ws.Range("E25").Value = '  -3.93%  '
$ws.Range("D26").Value = '170.27'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D28").Value = '0.137'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '1.43'
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("D30").Value = '19.59'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  +5.02%  '
$ws.Range("D32").Value = '0.121'
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").Value = '4.47'
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("D35").Value = '6.94'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").Value = '0.0615'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +7.29%  '
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = '4.93'
$ws.Range("E40").Value = '  +18.62%  '
$ws.Range("D41").Value = '101.99'
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("D42").Value = '0.0227'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").Value = '17.68'
$ws.Range("E43").Value = '  -2.45%  '
$ws.Range("D44").Value = '1.512.21'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("E45").Value = '  +1.03%  '
$ws.Range("D46").Value = '7.86'
$ws.Range("E46").Value = '  +0.38%  '
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("D48").Value = '0.0913'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("E50").Value = '  +31.93%  '
$ws.Range("E51").Value = '  +0.90%  '
